$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-18: price & volume updates
$ws.Range("D2").Value = "41.418.51"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "2.442.83"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.13"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.14"
$ws.Range("E6").Value = "  -2.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.543"
$ws.Range("E7").Value = "  -1.83%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -3.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.25"
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0831"
$ws.Range("E11").Value = "  -7.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.108"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").Value = "2.815.54"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.37"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").Value = "2.436.38"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.774"
$ws.Range("E17").Value = "  -1.94%  "
$ws.Range("D18").Value = "41.312.45"
$ws.Range("E18").Value = "  -0.72%  "

# Rows 19-20: ShibaInu and Uniswap swapped rank position
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0922"
$ws.Range("E19").Value = "  -4.17%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.26"
$ws.Range("E20").Value = "  -3.15%  "

# Rows 21-51: price & volume updates
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.11"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.08"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.62"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.11"
$ws.Range("E27").Value = "  -3.18%  "
$ws.Range("E28").Value = "  -3.26%  "
$ws.Range("E29").Value = "  -2.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.85"
$ws.Range("E30").Value = "  -3.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.53"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.27"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("E35").Value = "  -2.63%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.58"
$ws.Range("E37").Value = "  -4.75%  "
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("E40").Value = "  -2.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.90"
$ws.Range("E41").Value = "  -1.86%  "
$ws.Range("E42").Value = "  -7.20%  "
$ws.Range("D43").Value = "1.989.31"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("E44").Value = "  -3.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.06"
$ws.Range("E45").Value = "  -6.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.87"
$ws.Range("E46").Value = "  -3.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.46"
$ws.Range("E47").Value = "  +3.34%  "
$ws.Range("D48").Value = "2.678.02"
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "95.47"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.44"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.13"
$ws.Range("E51").Value = "  -0.66%  "
